$d = $word.ActiveDocument

# The document originally contains two custom-service field codes
# (" m:self.someCustomService() " and " m:self.someOtherCustomService() ")
# rendered as real Word fields (fldChar begin/instrText/fldChar end).
# The parser was switched to TokenIteratorFieldRewriterSplit, which expects
# plain template text like "{m:self.someCustomService()}" instead of a
# field. Rewrite each field's paragraph into plain runs, keeping the
# "self" run's orange theme color (accent6, shaded).

function Convert-CustomServiceField($suffix) {
    # Always operate on the first remaining field; deleting a field
    # re-indexes the collection, so we never need an explicit index.
    $f = $d.Fields.Item(1)
    $codeStart = $f.Code.Start

    # Locate the paragraph that owns this field so we know where to
    # reinsert plain text once the field is gone.
    $paraStart = $null
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($codeStart -ge $p.Range.Start -and $codeStart -lt $p.Range.End) {
            $paraStart = $p.Range.Start
        }
    }

    # Remove the field (begin mark / instrText runs / end mark) leaving an
    # empty paragraph behind.
    $f.Delete()

    $r = $d.Range($paraStart, $paraStart)

    # "{m:"
    $r.InsertAfter("{m:")
    $r.Collapse(0)

    # "self" in the original accent6/BF themed orange color
    $r.InsertAfter("self")
    $r.Font.TextColor.ObjectThemeColor = 9
    $r.Font.TextColor.TintAndShade = -0.25098039215686274
    $r.Collapse(0)

    # ".someCustomService()}" / ".someOtherCustomService()}"
    $r.InsertAfter($suffix)
}

Convert-CustomServiceField ".someCustomService()}"
Convert-CustomServiceField ".someOtherCustomService()}"
